$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch remain plain text,
# matching the original inlineStr content (values like "0.999" or
# "66.545.88" would otherwise be reinterpreted as numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.545.88'
$ws.Range("E2").Value = '  +2.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.490.91'
$ws.Range("E3").Value = '  +1.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.55'
$ws.Range("E5").Value = '  +1.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.16'
$ws.Range("E6").Value = '  +0.55%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.486.30'
$ws.Range("E8").Value = '  +1.49%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.591'
$ws.Range("E9").Value = '  +6.08%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  +5.33%  '

# Row 11
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.31'
$ws.Range("E11").Value = '  +0.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  +0.94%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.092.25'
$ws.Range("E13").Value = '  +1.71%  '

# Row 14
$ws.Range("E14").Value = '  -0.41%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.01'
$ws.Range("E15").Value = '  +2.77%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000178'
$ws.Range("E16").Value = '  +1.95%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.530.87'
$ws.Range("E17").Value = '  +2.16%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.472.67'
$ws.Range("E18").Value = '  +1.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("E19").Value = '  +0.88%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.94'
$ws.Range("E20").Value = '  +2.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.87'
$ws.Range("E21").Value = '  +2.31%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.88'
$ws.Range("E22").Value = '  -0.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.75'
$ws.Range("E23").Value = '  +2.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.530'
$ws.Range("E25").Value = '  +2.40%  '

# Row 26
$ws.Range("E26").Value = '  +2.86%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  +3.93%  '

# Row 28
$ws.Range("E28").Value = '  -0.36%  '

# Row 29
$ws.Range("E29").Value = '  +0.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.32'
$ws.Range("E30").Value = '  +2.53%  '

# Row 31
$ws.Range("E31").Value = '  +0.69%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("E32").Value = '  +1.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.65'
$ws.Range("E33").Value = '  +2.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.31'
$ws.Range("E34").Value = '  +2.19%  '

# Row 35
$ws.Range("E35").Value = '  +5.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.79'
$ws.Range("E36").Value = '  +2.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.895'
$ws.Range("E37").Value = '  +2.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.91'
$ws.Range("E38").Value = '  +2.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.80'
$ws.Range("E39").Value = '  +3.95%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.62'
$ws.Range("E40").Value = '  +4.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0736'
$ws.Range("E41").Value = '  +0.93%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.14'
$ws.Range("E42").Value = '  +0.26%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.58'
$ws.Range("E43").Value = '  +1.19%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.767.92'
$ws.Range("E44").Value = '  -0.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.79'
$ws.Range("E45").Value = '  -0.34%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.55'
$ws.Range("E46").Value = '  +2.43%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0308'
$ws.Range("E47").Value = '  +0.10%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '342.03'
$ws.Range("E48").Value = '  +3.38%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  +1.61%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.39'
$ws.Range("E50").Value = '  +7.66%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.852'
$ws.Range("E51").Value = '  +3.59%  '
